# Update the dSF column (F) values on Sheet1 to reflect the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -2
    3  = -1
    4  = 3
    5  = -3
    7  = -1
    8  = 4
    9  = -1
    10 = 3
    11 = -4
    12 = 3
    13 = 5
    14 = 2
    15 = 0
    16 = -2
    17 = 2
    19 = 2
    21 = -2
    24 = 2
    25 = 2
    26 = 2
    27 = -2
    28 = -4
    29 = -5
    30 = 6
    31 = -1
    33 = -1
    34 = -3
    35 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
